$wb = $excel.ActiveWorkbook

# "zh-cn" sheet - update handoff/handback correspond datetimes for the
# 1c57ed33... file row (row 2) to reflect the new handback run.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 02:47:20"
$wsZhCn.Range("H2").Value = "2016-03-23 02:47:52"

# "de-de" sheet - same update for the 1c57ed33... file row (row 2).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 02:47:25"
$wsDeDe.Range("H2").Value = "2016-03-23 02:47:58"
